# Apply the updated cryptocurrency price/volume figures (and the two
# reordered coin rows) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is: cell reference, new text value, and whether the value
# looks like a plain number and therefore needs to be forced back to
# text (every Price/Volume cell in this sheet is stored as text).
$updates = @(
    ,@("D2", "46.169.13", $false)
    ,@("E2", "  +1.72%  ", $false)
    ,@("D3", "2.606.13", $false)
    ,@("E3", "  +8.05%  ", $false)
    ,@("D4", "0.998", $true)
    ,@("E4", "  -0.11%  ", $false)
    ,@("D5", "308.04", $true)
    ,@("E5", "  +4.87%  ", $false)
    ,@("D6", "99.85", $true)
    ,@("E6", "  +5.89%  ", $false)
    ,@("D7", "0.603", $true)
    ,@("E7", "  +6.85%  ", $false)
    ,@("E8", "  +0.09%  ", $false)
    ,@("D9", "0.582", $true)
    ,@("E9", "  +16.38%  ", $false)
    ,@("D10", "39.51", $true)
    ,@("E10", "  +13.80%  ", $false)
    ,@("D11", "54.36", $true)
    ,@("E11", "  +1.73%  ", $false)
    ,@("D12", "0.0844", $true)
    ,@("E12", "  +8.49%  ", $false)
    ,@("D13", "8.22", $true)
    ,@("E13", "  +16.90%  ", $false)
    ,@("D14", "2.992.17", $false)
    ,@("E14", "  +7.85%  ", $false)
    ,@("D15", "0.106", $true)
    ,@("E15", "  +1.70%  ", $false)
    ,@("D16", "2.600.13", $false)
    ,@("E16", "  +7.65%  ", $false)
    ,@("D17", "0.922", $true)
    ,@("E17", "  +10.52%  ", $false)
    ,@("D18", "14.99", $true)
    ,@("E18", "  +6.69%  ", $false)
    ,@("D19", "46.346.03", $false)
    ,@("E19", "  +2.51%  ", $false)
    ,@("E20", "  +7.82%  ", $false)
    ,@("D21", "13.04", $true)
    ,@("E21", "  +5.82%  ", $false)
    ,@("D22", "6.73", $true)
    ,@("E22", "  +9.51%  ", $false)
    ,@("D23", "71.84", $true)
    ,@("E23", "  +7.33%  ", $false)
    ,@("D24", "272.69", $true)
    ,@("E24", "  +13.56%  ", $false)
    ,@("D25", "3.04", $true)
    ,@("E25", "  +9.83%  ", $false)
    ,@("D26", "30.39", $true)
    ,@("E26", "  +43.79%  ", $false)
    ,@("D27", "2.18", $true)
    ,@("E27", "  +12.63%  ", $false)
    ,@("E28", "  +0.39%  ", $false)
    ,@("E29", "  +0.40%  ", $false)
    ,@("D30", "10.60", $true)
    ,@("E30", "  +9.82%  ", $false)
    ,@("E31", "  +4.23%  ", $false)
    ,@("D32", "39.33", $true)
    ,@("E32", "  +1.70%  ", $false)
    ,@("D33", "6.22", $true)
    ,@("E33", "  +14.47%  ", $false)
    ,@("D34", "3.66", $true)
    ,@("E34", "  -2.49%  ", $false)
    ,@("D35", "2.84", $true)
    ,@("E35", "  +3.17%  ", $false)
    ,@("D36", "0.0842", $true)
    ,@("E36", "  +10.14%  ", $false)
    ,@("D37", "2.22", $true)
    ,@("E37", "  +11.85%  ", $false)
    ,@("D38", "150.00", $true)
    ,@("E38", "  +0.91%  ", $false)
    ,@("D39", "0.122", $true)
    ,@("E39", "  +7.83%  ", $false)
    ,@("E40", "  +6.21%  ", $false)
    ,@("D41", "23.25", $true)
    ,@("E41", "  +45.27%  ", $false)
    ,@("D42", "16.13", $true)
    ,@("E42", "  +8.82%  ", $false)
    ,@("B43", "NEARProtocol", $false)
    ,@("C43", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", $false)
    ,@("D43", "3.63", $true)
    ,@("E43", "  +13.16%  ", $false)
    ,@("B44", "VeChain", $false)
    ,@("C44", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", $false)
    ,@("D44", "0.0329", $true)
    ,@("E44", "  +10.81%  ", $false)
    ,@("D45", "4.11", $true)
    ,@("E45", "  +9.01%  ", $false)
    ,@("D46", "2.161.30", $false)
    ,@("E46", "  +8.42%  ", $false)
    ,@("D47", "0.997", $true)
    ,@("E47", "  -0.25%  ", $false)
    ,@("D48", "93.72", $true)
    ,@("E48", "  +5.68%  ", $false)
    ,@("D49", "9.64", $true)
    ,@("E49", "  +13.29%  ", $false)
    ,@("B50", "Aave", $false)
    ,@("C50", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", $false)
    ,@("D50", "109.36", $true)
    ,@("E50", "  +8.78%  ", $false)
    ,@("B51", "Stacks", $false)
    ,@("C51", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", $false)
    ,@("D51", "1.77", $true)
    ,@("E51", "  -1.14%  ", $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    $cell = $ws.Range($ref)
    if ($forceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
